# Add a new "Update" row (row 13) to the products sheet, mirroring the
# existing rows (same style, same cell typing) so the sheet's dimension
# grows from A1:G12 to A1:G13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / numeric columns: a direct Value assignment is fine since
# none of these strings look like dates/times to the auto-detect logic.
$ws.Range("A13").Value = "hcdihcid"
$ws.Range("B13").Value = "Flash USB"
$ws.Range("C13").Value = "128 GB"
$ws.Range("D13").Value = 90
$ws.Range("E13").Value = 2000
$ws.Range("G13").Value = "17:50:33"

# Column F holds a literal "YYYY-MM-DD" string (same as the other rows'
# Date column), not a real date value. Assigning that string straight to
# .Value would be auto-parsed into a date serial, which is not what the
# source workbook has (its Date column cells are inline/shared strings).
# Work around the auto-detection:
#   1) Clone F12's cell formatting (general, non-date) onto F13 so the
#      new cell doesn't pick up a freshly-minted number format.
#   2) Write the text through a formula (a quoted string literal is never
#      reinterpreted as a date).
#   3) Copy/PasteSpecial values-only over itself to flatten the formula
#      down to a plain literal string, leaving no formula behind.
$ws.Range("F12").Copy()
$ws.Range("F13").PasteSpecial(-4104)  # xlPasteAll (formatting + blank value)
$ws.Range("F13").Formula = "=""2024-09-23"""
$ws.Range("F13").Copy()
$ws.Range("F13").PasteSpecial(-4163)  # xlPasteValues (flatten formula to literal)

$excel.CutCopyMode = $false
